# "Generate Report for Handoff"
#
# The "8fc43afd-742d-4c56-8e46-f27024e5a5af" source file has just been
# handed off again, so its row now sorts first (alphabetically ahead of
# "2a294878-...") and keeps its "Handed back: in sync with en-US" status.
# The "2a294878-..." row moves to second place and flips to
# "Ready for handoff" because a fresh handoff was generated for it too
# (its handoff datetime is bumped to the same new handoff timestamp).
#
# Hyperlink target URLs (rId -> address) are untouched; only the visible
# cell text and the hyperlink's displayed text are updated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "8fc43afd-742d-4c56-8e46-f27024e5a5af.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"

$ws.Range("A3").Value = "2a294878-5d42-4456-b350-7b06c36b0e05.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "8fc43afd-742d-4c56-8e46-f27024e5a5af.md" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "2a294878-5d42-4456-b350-7b06c36b0e05.md" }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "8fc43afd-742d-4c56-8e46-f27024e5a5af.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "8fc43afd-742d-4c56-8e46-f27024e5a5af.6ae2fc2ffd2c31dc7888c39bd47754a5ec8b1b63.zh-cn.xlf"
$ws.Range("D2").Value = "2016-02-22 06:35:56"
$ws.Range("E2").Value = "8fc43afd-742d-4c56-8e46-f27024e5a5af.md"
$ws.Range("F2").Value = "8fc43afd-742d-4c56-8e46-f27024e5a5af.6ae2fc2ffd2c31dc7888c39bd47754a5ec8b1b63.zh-cn.xlf"
$ws.Range("G2").Value = "2016-02-22 06:34:33"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "2a294878-5d42-4456-b350-7b06c36b0e05.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "2a294878-5d42-4456-b350-7b06c36b0e05.7749c5a0c013bd4928f14ad39a73a9e25ca2e198.zh-cn.xlf"
$ws.Range("D3").Value = "2016-02-22 06:35:56"
$ws.Range("E3").Value = "2a294878-5d42-4456-b350-7b06c36b0e05.md"
$ws.Range("F3").Value = "2a294878-5d42-4456-b350-7b06c36b0e05.7749c5a0c013bd4928f14ad39a73a9e25ca2e198.zh-cn.xlf"
$ws.Range("G3").Value = "2016-02-22 06:34:33"
$ws.Range("H3").Value = "Include"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "8fc43afd-742d-4c56-8e46-f27024e5a5af.md" }
    elseif ($addr -eq '$C$2') { $h.TextToDisplay = "8fc43afd-742d-4c56-8e46-f27024e5a5af.6ae2fc2ffd2c31dc7888c39bd47754a5ec8b1b63.zh-cn.xlf" }
    elseif ($addr -eq '$E$2') { $h.TextToDisplay = "8fc43afd-742d-4c56-8e46-f27024e5a5af.md" }
    elseif ($addr -eq '$F$2') { $h.TextToDisplay = "8fc43afd-742d-4c56-8e46-f27024e5a5af.6ae2fc2ffd2c31dc7888c39bd47754a5ec8b1b63.zh-cn.xlf" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "2a294878-5d42-4456-b350-7b06c36b0e05.md" }
    elseif ($addr -eq '$C$3') { $h.TextToDisplay = "2a294878-5d42-4456-b350-7b06c36b0e05.7749c5a0c013bd4928f14ad39a73a9e25ca2e198.zh-cn.xlf" }
    elseif ($addr -eq '$E$3') { $h.TextToDisplay = "2a294878-5d42-4456-b350-7b06c36b0e05.md" }
    elseif ($addr -eq '$F$3') { $h.TextToDisplay = "2a294878-5d42-4456-b350-7b06c36b0e05.7749c5a0c013bd4928f14ad39a73a9e25ca2e198.zh-cn.xlf" }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "8fc43afd-742d-4c56-8e46-f27024e5a5af.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "8fc43afd-742d-4c56-8e46-f27024e5a5af.6ae2fc2ffd2c31dc7888c39bd47754a5ec8b1b63.de-de.xlf"
$ws.Range("D2").Value = "2016-02-22 06:36:11"
$ws.Range("E2").Value = "8fc43afd-742d-4c56-8e46-f27024e5a5af.md"
$ws.Range("F2").Value = "8fc43afd-742d-4c56-8e46-f27024e5a5af.6ae2fc2ffd2c31dc7888c39bd47754a5ec8b1b63.de-de.xlf"
$ws.Range("G2").Value = "2016-02-22 06:35:00"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "2a294878-5d42-4456-b350-7b06c36b0e05.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "2a294878-5d42-4456-b350-7b06c36b0e05.7749c5a0c013bd4928f14ad39a73a9e25ca2e198.de-de.xlf"
$ws.Range("D3").Value = "2016-02-22 06:36:11"
$ws.Range("E3").Value = "2a294878-5d42-4456-b350-7b06c36b0e05.md"
$ws.Range("F3").Value = "2a294878-5d42-4456-b350-7b06c36b0e05.7749c5a0c013bd4928f14ad39a73a9e25ca2e198.de-de.xlf"
$ws.Range("G3").Value = "2016-02-22 06:35:00"
$ws.Range("H3").Value = "Include"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "8fc43afd-742d-4c56-8e46-f27024e5a5af.md" }
    elseif ($addr -eq '$C$2') { $h.TextToDisplay = "8fc43afd-742d-4c56-8e46-f27024e5a5af.6ae2fc2ffd2c31dc7888c39bd47754a5ec8b1b63.de-de.xlf" }
    elseif ($addr -eq '$E$2') { $h.TextToDisplay = "8fc43afd-742d-4c56-8e46-f27024e5a5af.md" }
    elseif ($addr -eq '$F$2') { $h.TextToDisplay = "8fc43afd-742d-4c56-8e46-f27024e5a5af.6ae2fc2ffd2c31dc7888c39bd47754a5ec8b1b63.de-de.xlf" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "2a294878-5d42-4456-b350-7b06c36b0e05.md" }
    elseif ($addr -eq '$C$3') { $h.TextToDisplay = "2a294878-5d42-4456-b350-7b06c36b0e05.7749c5a0c013bd4928f14ad39a73a9e25ca2e198.de-de.xlf" }
    elseif ($addr -eq '$E$3') { $h.TextToDisplay = "2a294878-5d42-4456-b350-7b06c36b0e05.md" }
    elseif ($addr -eq '$F$3') { $h.TextToDisplay = "2a294878-5d42-4456-b350-7b06c36b0e05.7749c5a0c013bd4928f14ad39a73a9e25ca2e198.de-de.xlf" }
}

Write-Host "Done applying handoff report update."
